$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (target stored widths: A=14.42578125, B=14.7109375)
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334

# Update cell values (rows 1-4, columns A and B)
$ws.Range("A1").Value = -0.036061935545037706
$ws.Range("B1").Value = 0.036061934731249558

$ws.Range("A2").Value = 0.048853777235474516
$ws.Range("B2").Value = -0.048853778088672917

$ws.Range("A3").Value = -0.0018645753298280813
$ws.Range("B3").Value = 0.0018645745058398317

$ws.Range("A4").Value = 0.0013350510328165189
$ws.Range("B4").Value = -0.0013350519304597837
